$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 269, shifting the existing rows 269-288
# down to 273-292 (Excel copies formatting from the row above on insert,
# which correctly carries the date style on column D).
$ws.Range("A269:R272").EntireRow.Insert()

# Fill in the 4 newly inserted rows with the new week's data
# (Macroferia Regional de Talca - Melón, fecha 44585).

# Row 269: Calameño / Primera
$ws.Range("A269").Value = 5
$ws.Range("B269").Value = "Macroferia Regional de Talca"
$ws.Range("C269").Value = "Maule"
$ws.Range("D269").Value = 44585
$ws.Range("E269").Value = 7
$ws.Range("F269").Value = 100112027
$ws.Range("G269").Value = "Melón"
$ws.Range("H269").Value = "Calameño"
$ws.Range("I269").Value = "Primera"
$ws.Range("J269").Value = 4000
$ws.Range("K269").Value = 500
$ws.Range("L269").Value = 500
$ws.Range("M269").Value = 500
$ws.Range("N269").Value = "$/unidad"
$ws.Range("O269").Value = "Región del Maule"
$ws.Range("P269").Value = 500
$ws.Range("Q269").Value = 1
$ws.Range("R269").Value = "Hortaliza"

# Row 270: Calameño / Segunda
$ws.Range("A270").Value = 5
$ws.Range("B270").Value = "Macroferia Regional de Talca"
$ws.Range("C270").Value = "Maule"
$ws.Range("D270").Value = 44585
$ws.Range("E270").Value = 7
$ws.Range("F270").Value = 100112027
$ws.Range("G270").Value = "Melón"
$ws.Range("H270").Value = "Calameño"
$ws.Range("I270").Value = "Segunda"
$ws.Range("J270").Value = 4000
$ws.Range("K270").Value = 300
$ws.Range("L270").Value = 300
$ws.Range("M270").Value = 300
$ws.Range("N270").Value = "$/unidad"
$ws.Range("O270").Value = "Región del Maule"
$ws.Range("P270").Value = 300
$ws.Range("Q270").Value = 1
$ws.Range("R270").Value = "Hortaliza"

# Row 271: Tuna / Primera
$ws.Range("A271").Value = 5
$ws.Range("B271").Value = "Macroferia Regional de Talca"
$ws.Range("C271").Value = "Maule"
$ws.Range("D271").Value = 44585
$ws.Range("E271").Value = 7
$ws.Range("F271").Value = 100112027
$ws.Range("G271").Value = "Melón"
$ws.Range("H271").Value = "Tuna"
$ws.Range("I271").Value = "Primera"
$ws.Range("J271").Value = 3000
$ws.Range("K271").Value = 600
$ws.Range("L271").Value = 600
$ws.Range("M271").Value = 600
$ws.Range("N271").Value = "$/unidad"
$ws.Range("O271").Value = "Región del Maule"
$ws.Range("P271").Value = 600
$ws.Range("Q271").Value = 1
$ws.Range("R271").Value = "Hortaliza"

# Row 272: Tuna / Segunda
$ws.Range("A272").Value = 5
$ws.Range("B272").Value = "Macroferia Regional de Talca"
$ws.Range("C272").Value = "Maule"
$ws.Range("D272").Value = 44585
$ws.Range("E272").Value = 7
$ws.Range("F272").Value = 100112027
$ws.Range("G272").Value = "Melón"
$ws.Range("H272").Value = "Tuna"
$ws.Range("I272").Value = "Segunda"
$ws.Range("J272").Value = 3000
$ws.Range("K272").Value = 400
$ws.Range("L272").Value = 400
$ws.Range("M272").Value = 400
$ws.Range("N272").Value = "$/unidad"
$ws.Range("O272").Value = "Región del Maule"
$ws.Range("P272").Value = 400
$ws.Range("Q272").Value = 1
$ws.Range("R272").Value = "Hortaliza"
